# Update cryptocurrency price and 1h volume-change figures to the
# latest scrape values (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.536.83"
$ws.Range("E2").Value = "  -0.44%  "

$ws.Range("D3").Value = "2.628.89"
$ws.Range("E3").Value = "  -1.99%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.43%  "

$ws.Range("E7").Value = "  +0.08%  "

$ws.Range("E8").Value = "  -2.35%  "

$ws.Range("D9").Value = "2.628.97"
$ws.Range("E9").Value = "  -1.94%  "

$ws.Range("E10").Value = "  -2.62%  "

$ws.Range("E11").Value = "  +1.16%  "

$ws.Range("E12").Value = "  +1.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.70"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.85%  "

$ws.Range("D15").Value = "3.108.72"
$ws.Range("E15").Value = "  -1.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000182"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.91%  "

$ws.Range("D17").Value = "67.252.06"
$ws.Range("E17").Value = "  -0.66%  "

$ws.Range("D18").Value = "2.634.96"
$ws.Range("E18").Value = "  -1.45%  "

$ws.Range("E19").Value = "  +2.26%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.07"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.69%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "357.62"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.96%  "

$ws.Range("E22").Value = "  -1.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.67"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.48%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.46%  "

$ws.Range("E25").Value = "  +0.02%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "69.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.46%  "

$ws.Range("E28").Value = "  +1.37%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.62%  "

$ws.Range("E30").Value = "  -2.15%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "547.61"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.34%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.11%  "

$ws.Range("E34").Value = "  -2.10%  "

$ws.Range("E35").Value = "  +3.92%  "

$ws.Range("E36").Value = "  +0.08%  "

$ws.Range("E37").Value = "  -3.50%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "156.86"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.36%  "

$ws.Range("E39").Value = "  -2.81%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.366"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.10%  "

$ws.Range("E41").Value = "  -1.20%  "

$ws.Range("E42").Value = "  +1.84%  "

$ws.Range("E43").Value = "  -1.74%  "

$ws.Range("E44").Value = "  +0.06%  "

$ws.Range("E45").Value = "  -4.38%  "

$ws.Range("E46").Value = "  -0.70%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "152.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.53%  "

$ws.Range("E48").Value = "  -2.22%  "

$ws.Range("E49").Value = "  -1.62%  "

$ws.Range("E50").Value = "  -1.57%  "

$ws.Range("E51").Value = "  -1.12%  "
